$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The document currently has a single hidden "_GoBack" bookmark
#    sitting on the very last paragraph. The edit moves that bookmark
#    onto the newly-added "RETURN guess" paragraph, so first drop the
#    old one (while its name is still unambiguous) - it will be
#    recreated in the new spot below.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the first "END WHILE" paragraph (the one that ends the
#    "PROCEDURE user input" pseudocode block) and add a new paragraph
#    right after it containing "RETURN guess".
# ------------------------------------------------------------------
$found = $d.Content
[void]$found.Find.Execute("END WHILE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endWhilePara = $found.Paragraphs.Item(1)
$insertAt = $endWhilePara.Index

[void]$endWhilePara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($insertAt + 1)

# Build the new paragraph (matching indentation, text, and the
# relocated bookmark) as a raw WordprocessingML fragment so the
# bookmark start/end land exactly around the run, inside this single
# paragraph - precise positioning that simple Range/Bookmarks.Add
# calls can't guarantee.
$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t>RETURN guess</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$newPara.Range.InsertXML($newParaXml)
